$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Friend_Request_Management")
$ws.Activate()

# Update the values that decide whether to cancel or accept the friend request
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = 1

# Move the selection/active cell to B9 as recorded in the saved view state
$ws.Range("B9").Select()
